$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K" - strikeouts) regenerated values for rows 2-28
$kValues = @{
    2  = 14
    3  = 10
    4  = 1
    5  = 8
    6  = 7
    7  = 8
    8  = 6
    9  = 6
    10 = 6
    11 = 8
    12 = 4
    13 = 1
    14 = 9
    15 = 5
    16 = 8
    17 = 10
    18 = 4
    19 = 7
    20 = 10
    21 = 7
    22 = 9
    23 = 7
    24 = 7
    25 = 3
    26 = 5
    27 = 4
    28 = 5
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
